# Apply the data update described by the commit:
# "expermits todos no convexos menos el 5to"
# (update all the non-convex experiment values except the 5th one)
#
# This rewrites the numeric / expression values that live on the
# "Restricciones_del_follower", "Punto_modificado", "Vector_bf" and
# "Vector_BF" sheets. Every one of these cells is stored as TEXT in the
# workbook (even the ones that look like plain numbers), so we force a
# text number-format before writing the value and then restore the
# default "Normal" style so we don't leave a stray style behind.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Worksheet,
        [string]$Address,
        [string]$Value
    )
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $wsFollower "A2" "5.35 - 2x_1 + y_1 - y_2"
Set-TextValue $wsFollower "B2" "-2.8499999999999996"
Set-TextValue $wsFollower "D2" "0.73"
Set-TextValue $wsFollower "E2" "8.4"
Set-TextValue $wsFollower "F2" "8.8"

Set-TextValue $wsFollower "A3" "2.1499999999999932 + x_1 - 3x_2 + y_2"
Set-TextValue $wsFollower "B3" "-4.149999999999993"
Set-TextValue $wsFollower "D3" "0.59"
Set-TextValue $wsFollower "E3" "9.7"
Set-TextValue $wsFollower "F3" "0"

Set-TextValue $wsFollower "A4" "104.95 - y_1"
Set-TextValue $wsFollower "B4" "-104.95"
Set-TextValue $wsFollower "D4" "0.87"
Set-TextValue $wsFollower "F4" "2.8000000000000003"

Set-TextValue $wsFollower "A5" "-3.5999999999999996 - y_2"
Set-TextValue $wsFollower "B5" "-3.5999999999999996"
Set-TextValue $wsFollower "D5" "0.08"
Set-TextValue $wsFollower "E5" "0"
Set-TextValue $wsFollower "F5" "6.8999999999999995"

# ---------------------------------------------------------------------
# Punto_modificado
# ---------------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $wsPunto "A2" "53.35"
Set-TextValue $wsPunto "B2" "19.7"
Set-TextValue $wsPunto "C2" "104.95"
Set-TextValue $wsPunto "D2" "3.5999999999999996"

# ---------------------------------------------------------------------
# Vector_bf  (sheet name lookups are case-insensitive, and "Vector_bf" /
# "Vector_BF" only differ by case, so address these two by their fixed
# tab position instead of by name to avoid grabbing the wrong sheet)
# ---------------------------------------------------------------------
$wsVecbf = $wb.Worksheets.Item(5)

Set-TextValue $wsVecbf "A2" "4.14"
Set-TextValue $wsVecbf "A3" "-0.78"

# ---------------------------------------------------------------------
# Vector_BF
# ---------------------------------------------------------------------
$wsVecBF = $wb.Worksheets.Item(6)

Set-TextValue $wsVecBF "A2" "9.100000000000001"
Set-TextValue $wsVecBF "A3" "28.099999999999998"
Set-TextValue $wsVecBF "A4" "-8.9"
Set-TextValue $wsVecBF "A5" "-1.299999999999999"
